$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each inner array: Row, A, B, C, D, E, F, G, H, I, J, K, L, M, N, O, P, Q, R, S, T
$data = @(
    ,(2, "ECs", "Ucn2", "Il10rb", "ECs", 1, 0.3333333333333333, 0.2125036666666667, 0.637511, 0.1064749515943994, 0.1064749515943994, 3, 1, 20.49243633333333, 61.47730900000001, 0.1704171761399068, 0.1704171761399068, 4.354717859766556, 39.192460737899, 0.01814516058035081, 0.01814516058035081)
    ,(3, "ECs", "Ucn2", "Il10rb", "FAPs", 1, 0.3333333333333333, 0.2125036666666667, 0.637511, 0.1064749515943994, 0.1064749515943994, 3, 1, 10.84773033333333, 32.543191, 0.09021082417256794, 0.09021082417256795, 2.305182470844556, 20.746642237601, 0.009605193137065047, 0.009605193137065048)
    ,(4, "ECs", "Ucn2", "Il10rb", "M2", 1, 0.3333333333333333, 0.2125036666666667, 0.637511, 0.1064749515943994, 0.1064749515943994, 3, 1, 84.07831700000001, 252.234951, 0.6992038001079025, 0.6992038001079026, 17.86695064966234, 160.802555846961, 0.07444769077110902, 0.07444769077110905)
    ,(5, "ECs", "Ucn2", "Il10rb", "sCs", 1, 0.3333333333333333, 0.2125036666666667, 0.637511, 0.1064749515943994, 0.1064749515943994, 3, 1, 4.830172, 14.490516, 0.04016819957962274, 0.04016819957962274, 1.026429260630667, 9.237863345676001, 0.004276907105874505, 0.004276907105874506)
    ,(6, "FAPs", "Ucn2", "Il10rb", "ECs", 3, 1, 1.193141333333333, 3.579424, 0.597823405613129, 0.597823405613129, 3, 1, 20.49243633333333, 61.47730900000001, 0.1704171761399068, 0.1704171761399068, 24.45037281000178, 220.053355290016, 0.1018793766149315, 0.1018793766149315)
    ,(7, "FAPs", "Ucn2", "Il10rb", "FAPs", 3, 1, 1.193141333333333, 3.579424, 0.597823405613129, 0.597823405613129, 3, 1, 10.84773033333333, 32.543191, 0.09021082417256794, 0.09021082417256795, 12.94287543355378, 116.485878901984, 0.05393014213001174, 0.05393014213001175)
    ,(8, "FAPs", "Ucn2", "Il10rb", "M2", 3, 1, 1.193141333333333, 3.579424, 0.597823405613129, 0.597823405613129, 3, 1, 84.07831700000001, 252.234951, 0.6992038001079025, 0.6992038001079026, 100.3173152498027, 902.855837248224, 0.4180003969981477, 0.4180003969981478)
    ,(9, "FAPs", "Ucn2", "Il10rb", "sCs", 3, 1, 1.193141333333333, 3.579424, 0.597823405613129, 0.597823405613129, 3, 1, 4.830172, 14.490516, 0.04016819957962274, 0.04016819957962274, 5.763077860309334, 51.867700742784, 0.02401348987003792, 0.02401348987003792)
    ,(10, "M2", "Ucn2", "Il10rb", "ECs", 3, 1, 0.1501633333333333, 0.45049, 0.07523933068411523, 0.07523933068411523, 3, 1, 20.49243633333333, 61.47730900000001, 0.1704171761399068, 0.1704171761399068, 3.077212547934445, 27.69491293141, 0.01282207426984356, 0.01282207426984356)
    ,(11, "M2", "Ucn2", "Il10rb", "FAPs", 3, 1, 0.1501633333333333, 0.45049, 0.07523933068411523, 0.07523933068411523, 3, 1, 10.84773033333333, 32.543191, 0.09021082417256794, 0.09021082417256795, 1.628931345954445, 14.66038211359, 0.006787402031206415, 0.006787402031206416)
    ,(12, "M2", "Ucn2", "Il10rb", "M2", 3, 1, 0.1501633333333333, 0.45049, 0.07523933068411523, 0.07523933068411523, 3, 1, 84.07831700000001, 252.234951, 0.6992038001079025, 0.6992038001079026, 12.62548034177667, 113.62932307599, 0.05260762593190848, 0.05260762593190849)
    ,(13, "M2", "Ucn2", "Il10rb", "sCs", 3, 1, 0.1501633333333333, 0.45049, 0.07523933068411523, 0.07523933068411523, 3, 1, 4.830172, 14.490516, 0.04016819957962274, 0.04016819957962274, 0.7253147280933334, 6.52783255284, 0.003022228451156774, 0.003022228451156774)
    ,(14, "sCs", "Ucn2", "Il10rb", "ECs", 3, 1, 0.4400006666666667, 1.320002, 0.2204623121083564, 0.2204623121083564, 3, 1, 20.49243633333333, 61.47730900000001, 0.1704171761399068, 0.1704171761399068, 9.01668564829089, 81.15017083461801, 0.03757056467478088, 0.03757056467478088)
    ,(15, "sCs", "Ucn2", "Il10rb", "FAPs", 3, 1, 0.4400006666666667, 1.320002, 0.2204623121083564, 0.2204623121083564, 3, 1, 10.84773033333333, 32.543191, 0.09021082417256794, 0.09021082417256795, 4.773008578486889, 42.957077206382, 0.01988808687428473, 0.01988808687428474)
    ,(16, "sCs", "Ucn2", "Il10rb", "M2", 3, 1, 0.4400006666666667, 1.320002, 0.2204623121083564, 0.2204623121083564, 3, 1, 84.07831700000001, 252.234951, 0.6992038001079025, 0.6992038001079026, 36.99451553221135, 332.9506397899021, 0.1541480864067372, 0.1541480864067373)
    ,(17, "sCs", "Ucn2", "Il10rb", "sCs", 3, 1, 0.4400006666666667, 1.320002, 0.2204623121083564, 0.2204623121083564, 3, 1, 4.830172, 14.490516, 0.04016819957962274, 0.04016819957962274, 2.125278900114667, 19.127510101032, 0.008855574152553539, 0.008855574152553539)
)

foreach ($entry in $data) {
    $r = $entry[0]
    for ($c = 1; $c -le 20; $c++) {
        $ws.Cells.Item($r, $c).Value2 = $entry[$c]
    }
}